$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "SEPIO"
$ws.Range("C22").Value = "entity [BFO:0000001]"
$ws.Range("D22").Value = "research study [SEPIO:0000125]"
$ws.Range("E22").Value = "all"
